$d = $word.ActiveDocument
$t = $d.Tables(1)

# The table has 20 rows; the 5 data rows (with content) are rows 1, 5, 9, 13, 17.
# Each data row has 5 cells. We update the cell text in place (same formatting
# is preserved since we only change the run text), which reproduces the same
# end-state as the "insert 2 cells / drop 2 cells" reshuffle described in the
# commit diff for row 1 (net effect: same 5 values, same cell formatting).

# Row 1 (table row index 1)
$row = $t.Rows(1)
$row.Cells(1).Range.Text = "32÷3=10, 2"
$row.Cells(2).Range.Text = "34÷7=4, 6"
$row.Cells(3).Range.Text = "68÷7=9, 5"
$row.Cells(4).Range.Text = "66÷6=11, 0"
$row.Cells(5).Range.Text = "88÷2=44, 0"

# Row 5 (table row index 5)
$row = $t.Rows(5)
$row.Cells(1).Range.Text = "57÷2=28, 1"
$row.Cells(2).Range.Text = "17÷3=5, 2"
$row.Cells(3).Range.Text = "86÷2=43, 0"
$row.Cells(4).Range.Text = "15÷8=1, 7"
$row.Cells(5).Range.Text = "14÷8=1, 6"

# Row 9 (table row index 9)
$row = $t.Rows(9)
$row.Cells(1).Range.Text = "43÷2=21, 1"
$row.Cells(2).Range.Text = "51÷9=5, 6"
$row.Cells(3).Range.Text = "63÷8=7, 7"
$row.Cells(4).Range.Text = "66÷9=7, 3"
$row.Cells(5).Range.Text = "98÷6=16, 2"

# Row 13 (table row index 13)
$row = $t.Rows(13)
$row.Cells(1).Range.Text = "53÷3=17, 2"
$row.Cells(2).Range.Text = "38÷2=19, 0"
$row.Cells(3).Range.Text = "34÷4=8, 2"
$row.Cells(4).Range.Text = "97÷2=48, 1"
$row.Cells(5).Range.Text = "25÷2=12, 1"

# Row 17 (table row index 17)
$row = $t.Rows(17)
$row.Cells(1).Range.Text = "96÷6=16, 0"
$row.Cells(2).Range.Text = "88÷5=17, 3"
$row.Cells(3).Range.Text = "64÷7=9, 1"
$row.Cells(4).Range.Text = "19÷3=6, 1"
$row.Cells(5).Range.Text = "81÷7=11, 4"

Write-Output "All cell updates applied"
